$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 16692709
$ws.Cells.Item(17, 10).Value = 16692709
$ws.Cells.Item(17, 12).Value = 50078127
$ws.Cells.Item(17, 14).Value = -50078463

$ws.Cells.Item(97, 8).Value = 398.33334
$ws.Cells.Item(97, 10).Value = 398.33334
$ws.Cells.Item(97, 12).Value = 1195.00002
$ws.Cells.Item(97, 14).Value = -2187.00002

$ws.Cells.Item(116, 8).Value = 3692.5386
$ws.Cells.Item(116, 9).Value = 2633.2222
$ws.Cells.Item(116, 10).Value = 6076
$ws.Cells.Item(116, 11).Value = 2633.2222
$ws.Cells.Item(116, 12).Value = 6076
$ws.Cells.Item(116, 13).Value = 808.7777999999998
$ws.Cells.Item(116, 14).Value = -12960

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 15155093
$ws.Cells.Item(74, 9).Value = 17859788
$ws.Cells.Item(74, 10).Value = 8802.799999999999
$ws.Cells.Item(74, 11).Value = 17859788
$ws.Cells.Item(74, 12).Value = 8802.799999999999
$ws.Cells.Item(74, 13).Value = -17858914
$ws.Cells.Item(74, 14).Value = -10550.8

$ws.Cells.Item(77, 8).Value = 15155093
$ws.Cells.Item(77, 9).Value = 17859788
$ws.Cells.Item(77, 10).Value = 8802.799999999999
$ws.Cells.Item(77, 11).Value = 89298940
$ws.Cells.Item(77, 12).Value = 44014
$ws.Cells.Item(77, 13).Value = -89294572
$ws.Cells.Item(77, 14).Value = -52750

$ws.Cells.Item(139, 8).Value = 44746.668
$ws.Cells.Item(139, 10).Value = 44746.668
$ws.Cells.Item(139, 12).Value = 44746.668
$ws.Cells.Item(139, 14).Value = -55026.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1957
$ws.Cells.Item(20, 9).Value = 1986
$ws.Cells.Item(20, 10).Value = 1921.5555
$ws.Cells.Item(20, 11).Value = 1986
$ws.Cells.Item(20, 12).Value = 1921.5555
$ws.Cells.Item(20, 13).Value = -1739
$ws.Cells.Item(20, 14).Value = -2415.5555

$ws.Cells.Item(129, 8).Value = 48223
$ws.Cells.Item(129, 10).Value = 48223
$ws.Cells.Item(129, 12).Value = 48223
$ws.Cells.Item(129, 14).Value = -58223

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 9529263
$ws.Cells.Item(31, 9).Value = 5772.7915
$ws.Cells.Item(31, 10).Value = 30307788
$ws.Cells.Item(31, 11).Value = 5772.7915
$ws.Cells.Item(31, 12).Value = 30307788
$ws.Cells.Item(31, 13).Value = -5477.7915
$ws.Cells.Item(31, 14).Value = -30308378

$ws.Cells.Item(34, 8).Value = 9529263
$ws.Cells.Item(34, 9).Value = 5772.7915
$ws.Cells.Item(34, 10).Value = 30307788
$ws.Cells.Item(34, 11).Value = 5772.7915
$ws.Cells.Item(34, 12).Value = 30307788
$ws.Cells.Item(34, 13).Value = -5570.7915
$ws.Cells.Item(34, 14).Value = -30308192

$ws.Cells.Item(58, 8).Value = 1940.2
$ws.Cells.Item(58, 9).Value = 741.41174
$ws.Cells.Item(58, 11).Value = 741.41174
$ws.Cells.Item(58, 13).Value = -538.41174

$ws.Cells.Item(136, 8).Value = 1940.2
$ws.Cells.Item(136, 9).Value = 741.41174
$ws.Cells.Item(136, 11).Value = 2224.23522
$ws.Cells.Item(136, 13).Value = 325.76478

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 679.875
$ws.Cells.Item(121, 9).Value = 396.66666
$ws.Cells.Item(121, 10).Value = 849.8
$ws.Cells.Item(121, 11).Value = 1189.99998
$ws.Cells.Item(121, 12).Value = 2549.4
$ws.Cells.Item(121, 13).Value = 120.0000199999999
$ws.Cells.Item(121, 14).Value = -5169.4

$ws.Cells.Item(129, 8).Value = 3568.7
$ws.Cells.Item(129, 9).Value = 1201
$ws.Cells.Item(129, 10).Value = 5147.1665
$ws.Cells.Item(129, 11).Value = 3603
$ws.Cells.Item(129, 12).Value = 15441.4995
$ws.Cells.Item(129, 13).Value = 1397
$ws.Cells.Item(129, 14).Value = -25441.4995

$ws.Cells.Item(131, 8).Value = 1148.4286
$ws.Cells.Item(131, 10).Value = 1211.9474
$ws.Cells.Item(131, 12).Value = 3635.8422
$ws.Cells.Item(131, 14).Value = -13715.8422

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3491.1482
$ws.Cells.Item(126, 9).Value = 2274.75
$ws.Cells.Item(126, 10).Value = 4464.2666
$ws.Cells.Item(126, 11).Value = 6824.25
$ws.Cells.Item(126, 12).Value = 13392.7998
$ws.Cells.Item(126, 13).Value = -4354.25
$ws.Cells.Item(126, 14).Value = -18332.7998

$ws.Cells.Item(135, 8).Value = 41199.75
$ws.Cells.Item(135, 10).Value = 41199.75
$ws.Cells.Item(135, 12).Value = 41199.75
$ws.Cells.Item(135, 14).Value = -51339.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4985.943
$ws.Cells.Item(7, 9).Value = 4442.5264
$ws.Cells.Item(7, 10).Value = 5631.25
$ws.Cells.Item(7, 11).Value = 4442.5264
$ws.Cells.Item(7, 12).Value = 5631.25
$ws.Cells.Item(7, 13).Value = -4330.5264
$ws.Cells.Item(7, 14).Value = -5855.25

$ws.Cells.Item(40, 8).Value = 4257.8184
$ws.Cells.Item(40, 9).Value = 5741.5835
$ws.Cells.Item(40, 10).Value = 3409.9524
$ws.Cells.Item(40, 11).Value = 5741.5835
$ws.Cells.Item(40, 12).Value = 3409.9524
$ws.Cells.Item(40, 13).Value = -5605.5835
$ws.Cells.Item(40, 14).Value = -3681.9524

$ws.Cells.Item(82, 8).Value = 2750
$ws.Cells.Item(82, 9).Value = 2333.3333
$ws.Cells.Item(82, 10).Value = 3000
$ws.Cells.Item(82, 11).Value = 2333.3333
$ws.Cells.Item(82, 12).Value = 3000
$ws.Cells.Item(82, 13).Value = -1972.3333
$ws.Cells.Item(82, 14).Value = -3722

$ws.Cells.Item(85, 8).Value = 2750
$ws.Cells.Item(85, 9).Value = 2333.3333
$ws.Cells.Item(85, 10).Value = 3000
$ws.Cells.Item(85, 11).Value = 2333.3333
$ws.Cells.Item(85, 12).Value = 3000
$ws.Cells.Item(85, 13).Value = -1085.3333
$ws.Cells.Item(85, 14).Value = -5496

$ws.Cells.Item(122, 8).Value = 4907.357
$ws.Cells.Item(122, 9).Value = 4915.96
$ws.Cells.Item(122, 10).Value = 4894.706
$ws.Cells.Item(122, 11).Value = 14747.88
$ws.Cells.Item(122, 12).Value = 14684.118
$ws.Cells.Item(122, 13).Value = -12297.88
$ws.Cells.Item(122, 14).Value = -19584.118

$ws.Cells.Item(126, 8).Value = 4985.943
$ws.Cells.Item(126, 9).Value = 4442.5264
$ws.Cells.Item(126, 10).Value = 5631.25
$ws.Cells.Item(126, 11).Value = 13327.5792
$ws.Cells.Item(126, 12).Value = 16893.75
$ws.Cells.Item(126, 13).Value = -10857.5792
$ws.Cells.Item(126, 14).Value = -21833.75

$ws.Cells.Item(132, 8).Value = 8938869
$ws.Cells.Item(132, 9).Value = 8360.451999999999
$ws.Cells.Item(132, 10).Value = 35730390
$ws.Cells.Item(132, 11).Value = 25081.356
$ws.Cells.Item(132, 12).Value = 107191170
$ws.Cells.Item(132, 13).Value = -22551.356
$ws.Cells.Item(132, 14).Value = -107196230

$ws.Cells.Item(136, 8).Value = 18527448
$ws.Cells.Item(136, 9).Value = 20834960
$ws.Cells.Item(136, 11).Value = 62504880
$ws.Cells.Item(136, 13).Value = -62502330

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2989.2173
$ws.Cells.Item(126, 9).Value = 2092
$ws.Cells.Item(126, 10).Value = 7251
$ws.Cells.Item(126, 11).Value = 6276
$ws.Cells.Item(126, 12).Value = 21753
$ws.Cells.Item(126, 13).Value = -3806
$ws.Cells.Item(126, 14).Value = -26693

$ws.Cells.Item(136, 8).Value = 762.48486
$ws.Cells.Item(136, 9).Value = 755.0625
$ws.Cells.Item(136, 10).Value = 1000
$ws.Cells.Item(136, 11).Value = 2265.1875
$ws.Cells.Item(136, 12).Value = 3000
$ws.Cells.Item(136, 13).Value = 284.8125
$ws.Cells.Item(136, 14).Value = -8100

Write-Host "All edits applied"